$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("特变电工", "金风科技", "海格通信"),
    @("中国西电", "特变电工", "岩山科技"),
    @("金风科技", "岩山科技", "特变电工"),
    @("五洲新春", "中国西电", "金风科技"),
    @("岩山科技", "海格通信", "三变科技"),
    @("锋龙股份", "保变电气", "中国西电"),
    @("海格通信", "汉缆股份", "平潭发展"),
    @("汉缆股份", "蓝色光标", "航天发展"),
    @("保变电气", "长电科技", "保变电气"),
    @("长电科技", "五洲新春", "神剑股份"),
    @("航天发展", "航天电子", "五洲新春"),
    @("三变科技", "航天发展", "锋龙股份"),
    @("森源电气", "三花智控", "汉缆股份"),
    @("中国卫星", "森源电气", "华夏幸福"),
    @("三花智控", "双杰电气", "国晟科技"),
    @("蓝色光标", "三变科技", "利欧股份"),
    @("双杰电气", "中国卫星", "三花智控"),
    @("航天电子", "华胜天成", "华胜天成"),
    @("华胜天成", "锋龙股份", "露笑科技"),
    @("超捷股份", "康强电子", "新联电子")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
